# Auto-generated edit script to update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.977.90"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").Value = "1.637.52"
$ws.Range("E3").Value = "  -0.36%  "

# Row 4
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.83%  "

# Row 5
$ws.Range("D5").Value = "'214.63"
$ws.Range("E5").Value = "  -0.50%  "

# Row 6
$ws.Range("D6").Value = "'0.5091"
$ws.Range("E6").Value = "  +0.70%  "

# Row 7
$ws.Range("E7").Value = "  -0.49%  "

# Row 8
$ws.Range("D8").Value = "'0.2567"

# Row 9
$ws.Range("D9").Value = "'0.06355"
$ws.Range("E9").Value = "  -0.75%  "

# Row 10
$ws.Range("D10").Value = "'19.67"
$ws.Range("E10").Value = "  +0.27%  "

# Row 11
$ws.Range("D11").Value = "'0.07773"
$ws.Range("E11").Value = "  -0.32%  "

# Row 12
$ws.Range("D12").Value = "'4.269"
$ws.Range("E12").Value = "  -0.47%  "

# Row 13
$ws.Range("D13").Value = "1.635.90"
$ws.Range("E13").Value = "  -0.75%  "

# Row 14
$ws.Range("D14").Value = "'0.5438"
$ws.Range("E14").Value = "  -0.10%  "

# Row 15
$ws.Range("D15").Value = "0.0₅7714"
$ws.Range("E15").Value = "  -2.04%  "

# Row 16
$ws.Range("D16").Value = "'63.98"
$ws.Range("E16").Value = "  -1.51%  "

# Row 17
$ws.Range("D17").Value = "25.983.17"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  -0.66%  "

# Row 19
$ws.Range("D19").Value = "'199.17"
$ws.Range("E19").Value = "  +0.62%  "

# Row 20
$ws.Range("D20").Value = "'4.418"
$ws.Range("E20").Value = "  +0.00%  "

# Row 21
$ws.Range("D21").Value = "'9.908"
$ws.Range("E21").Value = "  -0.69%  "

# Row 22
$ws.Range("D22").Value = "'6.043"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  -0.49%  "

# Row 24
$ws.Range("E24").Value = "  +1.13%  "

# Row 25
$ws.Range("D25").Value = "'141.25"
$ws.Range("E25").Value = "  +0.34%  "

# Row 26
$ws.Range("D26").Value = "'0.1204"
$ws.Range("E26").Value = "  +5.20%  "

# Row 27
$ws.Range("D27").Value = "'6.830"
$ws.Range("E27").Value = "  -0.80%  "

# Row 28
$ws.Range("D28").Value = "'15.63"
$ws.Range("E28").Value = "  -0.81%  "

# Row 29
$ws.Range("D29").Value = "'1.233"
$ws.Range("E29").Value = "  -1.03%  "

# Row 30
$ws.Range("D30").Value = "'0.04897"
$ws.Range("E30").Value = "  -2.42%  "

# Row 31
$ws.Range("D31").Value = "'3.262"
$ws.Range("E31").Value = "  -0.28%  "

# Row 32
$ws.Range("D32").Value = "'3.174"

# Row 33
$ws.Range("D33").Value = "'1.530"
$ws.Range("E33").Value = "  -0.37%  "

# Row 34
$ws.Range("D34").Value = "'2.371"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("D35").Value = "'0.9085"
$ws.Range("E35").Value = "  +1.49%  "

# Row 36
$ws.Range("D36").Value = "'2.586"
$ws.Range("E36").Value = "  -1.23%  "

# Row 37
$ws.Range("D37").Value = "1.127.79"
$ws.Range("E37").Value = "  -1.51%  "

# Row 38
$ws.Range("D38").Value = "'0.5464"
$ws.Range("E38").Value = "  -1.65%  "

# Row 39
$ws.Range("D39").Value = "'0.01562"
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("E40").Value = "  -0.53%  "

# Row 41
$ws.Range("E41").Value = "  -1.55%  "

# Row 42
$ws.Range("D42").Value = "'0.8115"
$ws.Range("E42").Value = "  -1.58%  "

# Row 43
$ws.Range("D43").Value = "0.0₈124"
$ws.Range("E43").Value = "  +1.99%  "

# Row 44
$ws.Range("D44").Value = "'98.93"
$ws.Range("E44").Value = "  -1.02%  "

# Row 45
$ws.Range("D45").Value = "'5.435"
$ws.Range("E45").Value = "  -4.58%  "

# Row 46
$ws.Range("D46").Value = "1.775.98"

# Row 47
$ws.Range("D47").Value = "'0.4521"
$ws.Range("E47").Value = "  -0.49%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'54.94"
$ws.Range("E48").Value = "  -0.95%  "

# Row 49
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "'0.9980"
$ws.Range("E49").Value = "  -0.82%  "

# Row 50
$ws.Range("D50").Value = "'0.05114"
$ws.Range("E50").Value = "  +0.75%  "

# Row 51
$ws.Range("E51").Value = "  -0.08%  "
